$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row: "<Column>_old" -> "<Column>_FV2310" (cols A-J)
#    and "<Column>_new" -> "<Column>_FV2404" (cols L-U). Column K ("diff")
#    is left untouched.
$fv2310Headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)
$fv2404Headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $fv2310Headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2310Headers[$i]
}

# column 11 ("K1") is "diff" - unchanged

for ($i = 0; $i -lt $fv2404Headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2404Headers[$i]
}

# 2. Create an Excel Table (ListObject) covering the used range
$rng = $ws.Range("A1:U68")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = $null

# 3. Freeze the header row (top row split) and set selection on bottom-left pane
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
